# The sheet originally has per-year column headers stored as text labels
# like "1960 [YR1960]" .. "2019 [YR2019]" in row 1 (columns E..BL). The edit
# converts those to plain numeric years (1960 .. 2019) and left-aligns them
# (to match how "2020 [YR2020]" in BM1, which stays textual, reads in the
# header row). This is what enables the year axis to be used numerically,
# e.g. for a line graph / correlation against the data rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 5   # column E  -> 1960
$lastCol  = 64  # column BL -> 2019
$startYear = 1960
$xlLeft = -4131 # xlHAlignLeft

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $year = $startYear + ($col - $firstCol)
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $year
    $cell.HorizontalAlignment = $xlLeft
}

# Reflect the header row as the active selection, matching the edited file.
$ws.Range("E1:BL1").Select() | Out-Null
